$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 158, shifting rows 158:208 down to 159:209
$ws.Rows("158:158").Insert()

# Populate the newly inserted row 158 with the new price-report record
$ws.Range("A158").Value = 3
$ws.Range("B158").Value = 'Femacal de La Calera'
$ws.Range("C158").Value = 'Coquimbo'
$ws.Range("D158").Value = 44468
$ws.Range("E158").Value = 5
$ws.Range("F158").Value = 100112040
$ws.Range("G158").Value = 'Cilantro'
$ws.Range("H158").Value = 'Sin especificar'
$ws.Range("I158").Value = 'Primera'
$ws.Range("J158").Value = 300
$ws.Range("K158").Value = 2500
$ws.Range("L158").Value = 3000
$ws.Range("M158").Value = 2800
$ws.Range("N158").Value = '$/docena de atados (3 kilos)'
$ws.Range("O158").Value = 'Provincia de Quillota'
$ws.Range("P158").Value = 933
$ws.Range("Q158").Value = 3
$ws.Range("R158").Value = 'Hortaliza'
